$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Overview sheet: the status summary text changes from
# "Ready for handoff" to "Handoff transform failed" (B2 and C2 both
# reference the same shared string, so both columns update).
# ------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B2").Value = "Handoff transform failed"
$wsOverview.Range("C2").Value = "Handoff transform failed"

# ------------------------------------------------------------------
# Helper: apply the "handoff transform failed" reset to a language
# sheet (zh-cn / de-de). The handoff attempt is rolled back: the
# "Latest Handoff File" cell (C2) is cleared out completely (along
# with its hyperlink), "Latest Handoff Datetime" / "Latest Handback
# DateTime" go back to the zero date, and "Handoff Reason" becomes
# "Ignored" instead of "Include".
# ------------------------------------------------------------------
function Reset-HandoffSheet($ws, $mdUrl, $configUrl) {
    $ws.Range("B2").Value = "Handoff transform failed"

    # Drop the "Latest Handoff File" cell and its hyperlink entirely.
    $ws.Range("C2").Clear()

    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Ignored"

    # Deleting a single hyperlink item isn't reliable in this host, but
    # clearing the whole collection and re-adding the two that must
    # survive (A2, A3) reproduces the target hyperlinks/rels exactly.
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $mdUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, "ae013254-c540-4bff-a548-43c6ef4ab4af.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $configUrl, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, ".localization-config")
}

$mdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/75a9f97e188cbc2b2874046728f955f15de415b7/e2e/ae013254-c540-4bff-a548-43c6ef4ab4af.md"
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/75a9f97e188cbc2b2874046728f955f15de415b7/.localization-config"

# --- zh-cn sheet ---
$wsZh = $wb.Worksheets.Item("zh-cn")
Reset-HandoffSheet $wsZh $mdUrl $configUrl

# --- de-de sheet ---
$wsDe = $wb.Worksheets.Item("de-de")
Reset-HandoffSheet $wsDe $mdUrl $configUrl
